$wb = $excel.ActiveWorkbook

$oldGuid = "8a4a9857-0636-4a64-9576-3bad229a52b0"
$newGuid = "05e3eb0f-55c7-48b4-8b4f-120ac4583c02"
$oldZhHash = "6f133a0c096f6747adac6c4488f9670eb23a3487"
$newZhHash = "a4093e37184df47caf5446ed8a060af3cbadc1b8"
$oldDeHash = "6f133a0c096f6747adac6c4488f9670eb23a3487"
$newDeHash = "a4093e37184df47caf5446ed8a060af3cbadc1b8"

$newTopDate = "2016-08-25 02:57:12"
$newZhDate  = "2016-08-25 02:57:04"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newTopDate

# --- Sheet "zh-cn" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $newZhDate

# --- Sheet "de-de" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newDeHash.de-de.xlf"
$wsDe.Range("H2").Value = $newTopDate
